$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")
$ws.Activate()
$ws.Range("A4").Select()
$ws.Range("A4").Value = "Death's End"
$ws.Range("B4").Value = "Cixin Liu"
$ws.Range("C4").Value = 45338
$ws.Range("D4").Value = 45365
$ws.Range("E4").Value = "***"
$ws.Range("F4").Value = "each book had a wild ending. Great series. Made me more scared of the universe whilst wishing lightspeed travel become a reality for humanity"
